$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format first, otherwise Excel auto-converts the typed string into
# a numeric value (losing the original text formatting, e.g. trailing zeros).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D16", "D19", "D20", "D21", "D23", "D24", "D25", "D27", "D28", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated crypto price / volume data scraped this run.
$ws.Range("D2").Value = '43.818.07'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '2.313.93'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '114.37'
$ws.Range("E5").Value = '  +20.22%  '
$ws.Range("D6").Value = '270.15'
$ws.Range("E6").Value = '  +1.14%  '
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("D10").Value = '47.48'
$ws.Range("E10").Value = '  +7.38%  '
$ws.Range("D11").Value = '0.0945'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").Value = '8.88'
$ws.Range("E12").Value = '  +15.20%  '
$ws.Range("E13").Value = '  +2.09%  '
$ws.Range("D14").Value = '15.72'
$ws.Range("E14").Value = '  +3.68%  '
$ws.Range("D15").Value = '2.664.72'
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").Value = '0.860'
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").Value = '2.312.69'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '43.722.95'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = '0.0000110'
$ws.Range("E19").Value = '  +3.32%  '
$ws.Range("D20").Value = '6.68'
$ws.Range("E20").Value = '  +8.42%  '
$ws.Range("D21").Value = '72.75'
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("E22").Value = '  +6.58%  '
$ws.Range("D23").Value = '234.03'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").Value = '2.92'
$ws.Range("E24").Value = '  +17.09%  '
$ws.Range("D25").Value = '9.48'
$ws.Range("E25").Value = '  +6.23%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '11.46'
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("D28").Value = '42.31'
$ws.Range("E28").Value = '  +8.85%  '
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").Value = '177.70'
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("D32").Value = '21.94'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '0.0932'
$ws.Range("E33").Value = '  +5.61%  '
$ws.Range("D34").Value = '5.58'
$ws.Range("E34").Value = '  +4.34%  '
$ws.Range("E35").Value = '  +0.83%  '
$ws.Range("E36").Value = '  +6.67%  '
$ws.Range("D37").Value = '0.111'
$ws.Range("E37").Value = '  +3.51%  '
$ws.Range("D38").Value = '3.95'
$ws.Range("E38").Value = '  +21.15%  '
$ws.Range("D39").Value = '0.0357'
$ws.Range("E39").Value = '  +0.61%  '
$ws.Range("D40").Value = '0.244'
$ws.Range("E40").Value = '  +3.58%  '
$ws.Range("E41").Value = '  +1.18%  '
$ws.Range("D42").Value = '71.06'
$ws.Range("E42").Value = '  +13.91%  '
$ws.Range("D43").Value = '12.96'
$ws.Range("E43").Value = '  +9.40%  '
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '6.01'
$ws.Range("E44").Value = '  +15.35%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("E46").Value = '  +3.19%  '
$ws.Range("D47").Value = '8.83'
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("E49").Value = '  +11.47%  '
$ws.Range("D50").Value = '100.50'
$ws.Range("E50").Value = '  +2.21%  '
$ws.Range("D51").Value = '1.23'
$ws.Range("E51").Value = '  +3.60%  '
